$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from 2023-10-25 (45224) to 2023-11-03 (45233)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
